# Auto-generated edit script applying the Lich_Profits diff to the
# 8 leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 42
$ws.Range("H42").Value = 128
$ws.Range("I42").Value = 128
$ws.Range("K42").Value = 384
$ws.Range("M42").Value = -154
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 62
$ws.Range("H62").Value = 100022000
$ws.Range("I62").Value = 125002500
$ws.Range("K62").Value = 125002500
$ws.Range("M62").Value = -125001876
# Row 65
$ws.Range("H65").Value = 100022000
$ws.Range("I65").Value = 125002500
$ws.Range("K65").Value = 625012500
$ws.Range("M65").Value = -625009380
# Row 100
$ws.Range("H100").Value = 4205.3335
$ws.Range("I100").Value = 1696.4
$ws.Range("J100").Value = 16750
$ws.Range("K100").Value = 1696.4
$ws.Range("L100").Value = 16750
$ws.Range("M100").Value = -1155.4
$ws.Range("N100").Value = -17832
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 132
$ws.Range("H132").Value = 2160.4443
$ws.Range("I132").Value = 2169.8823
$ws.Range("K132").Value = 6509.646900000001
$ws.Range("M132").Value = -3979.646900000001
# Row 135
$ws.Range("H135").Value = 3326.4
$ws.Range("I135").Value = 3224.25
$ws.Range("K135").Value = 29018.25
$ws.Range("M135").Value = -26483.25
# Row 138
$ws.Range("H138").Value = 2751.21
$ws.Range("I138").Value = 1359.2069
$ws.Range("J138").Value = 3319.7747
$ws.Range("K138").Value = 4077.620699999999
$ws.Range("L138").Value = 9959.3241
$ws.Range("M138").Value = 1062.379300000001
$ws.Range("N138").Value = -20239.3241

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 95
$ws.Range("H95").Value = 22633.334
$ws.Range("J95").Value = 22633.334
$ws.Range("L95").Value = 22633.334
$ws.Range("N95").Value = -28125.334
# Row 101
$ws.Range("H101").Value = 29500
$ws.Range("J101").Value = 29500
$ws.Range("L101").Value = 29500
$ws.Range("N101").Value = -35990
# Row 122
$ws.Range("H122").Value = 5132.25
$ws.Range("I122").Value = 4113.7617
$ws.Range("K122").Value = 12341.2851
$ws.Range("M122").Value = -9891.285100000001

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 86
$ws.Range("H86").Value = 2131.9285
$ws.Range("I86").Value = 2153.8333
$ws.Range("J86").Value = 2000.5
$ws.Range("K86").Value = 2153.8333
$ws.Range("L86").Value = 2000.5
$ws.Range("M86").Value = -1030.8333
$ws.Range("N86").Value = -4246.5
# Row 89
$ws.Range("H89").Value = 2131.9285
$ws.Range("I89").Value = 2153.8333
$ws.Range("J89").Value = 2000.5
$ws.Range("K89").Value = 10769.1665
$ws.Range("L89").Value = 10002.5
$ws.Range("M89").Value = -5153.166499999999
$ws.Range("N89").Value = -21234.5
# Row 107
$ws.Range("H107").Value = 1615.6364
$ws.Range("I107").Value = 1577.2
$ws.Range("K107").Value = 1577.2
$ws.Range("M107").Value = 342.8

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 4
$ws.Range("H4").Value = 500002500
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888
# Row 7
$ws.Range("H7").Value = 1661
$ws.Range("J7").Value = 1922
$ws.Range("L7").Value = 1922
$ws.Range("N7").Value = -2148
# Row 31
$ws.Range("H31").Value = 402176.7
$ws.Range("J31").Value = 3435.3635
$ws.Range("L31").Value = 3435.3635
$ws.Range("N31").Value = -4025.3635
# Row 34
$ws.Range("H34").Value = 402176.7
$ws.Range("J34").Value = 3435.3635
$ws.Range("L34").Value = 3435.3635
$ws.Range("N34").Value = -3839.3635
# Row 42
$ws.Range("H42").Value = 16099.444
$ws.Range("J42").Value = 18285
$ws.Range("L42").Value = 18285
$ws.Range("N42").Value = -19471
# Row 43
$ws.Range("H43").Value = 12337.125
$ws.Range("J43").Value = 12337.125
$ws.Range("L43").Value = 12337.125
$ws.Range("N43").Value = -12705.125
# Row 101
$ws.Range("H101").Value = 12337.125
$ws.Range("J101").Value = 12337.125
$ws.Range("L101").Value = 12337.125
$ws.Range("N101").Value = -18827.125
# Row 105
$ws.Range("H105").Value = 4497.3022
$ws.Range("I105").Value = 1873.3636
$ws.Range("J105").Value = 7246.1904
$ws.Range("K105").Value = 1873.3636
$ws.Range("L105").Value = 7246.1904
$ws.Range("M105").Value = -126.3635999999999
$ws.Range("N105").Value = -10740.1904

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 37
$ws.Range("H37").Value = 91014730
$ws.Range("J37").Value = 91014730
$ws.Range("L37").Value = 273044190
$ws.Range("N37").Value = -273044414
# Row 86
$ws.Range("H86").Value = 938.8570999999999
$ws.Range("I86").Value = 938.8570999999999
$ws.Range("K86").Value = 2816.5713
$ws.Range("M86").Value = -1630.5713
# Row 89
$ws.Range("H89").Value = 938.8570999999999
$ws.Range("I89").Value = 938.8570999999999
$ws.Range("K89").Value = 8449.713899999999
$ws.Range("M89").Value = -2521.713899999999
# Row 92
$ws.Range("H92").Value = 624.25
$ws.Range("I92").Value = 749
$ws.Range("K92").Value = 2247
$ws.Range("M92").Value = -999
# Row 107
$ws.Range("H107").Value = 483.83334
$ws.Range("I107").Value = 472.7857
$ws.Range("J107").Value = 499.3
$ws.Range("K107").Value = 1418.3571
$ws.Range("L107").Value = 1497.9
$ws.Range("M107").Value = 501.6428999999998
$ws.Range("N107").Value = -5337.9
# Row 139
$ws.Range("H139").Value = 3459.3333
$ws.Range("I139").Value = 3019.889
$ws.Range("J139").Value = 4777.6665
$ws.Range("K139").Value = 9059.667000000001
$ws.Range("L139").Value = 14332.9995
$ws.Range("M139").Value = -3919.667000000001
$ws.Range("N139").Value = -24612.9995

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 15
$ws.Range("H15").Value = 29990
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 29990
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 29990
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -30566
# Row 81
$ws.Range("H81").Value = 29990
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 29990
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 29990
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -31986
# Row 84
$ws.Range("H84").Value = 29990
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 29990
$ws.Range("K84").Value = 89997
$ws.Range("L84").Value = 45000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -99954
# Row 92
$ws.Range("H92").Value = 5750.3335
$ws.Range("J92").Value = 5750.3335
$ws.Range("L92").Value = 5750.3335
$ws.Range("N92").Value = -9494.333500000001
# Row 102
$ws.Range("H102").Value = 3906.8572
$ws.Range("I102").Value = 3943.8823
$ws.Range("K102").Value = 3943.8823
$ws.Range("M102").Value = -2321.8823
# Row 126
$ws.Range("H126").Value = 9988
$ws.Range("I126").Value = 16673.75
$ws.Range("J126").Value = 3302.25
$ws.Range("K126").Value = 50021.25
$ws.Range("L126").Value = 9906.75
$ws.Range("M126").Value = -47551.25
$ws.Range("N126").Value = -14846.75
# Row 132
$ws.Range("H132").Value = 41145.348
$ws.Range("I132").Value = 44407.543
$ws.Range("K132").Value = 133222.629
$ws.Range("M132").Value = -130692.629

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 16
$ws.Range("H16").Value = 78231550
$ws.Range("I16").Value = 78231550
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 78231550
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -78231380
$ws.Range("N16").ClearContents()
# Row 22
$ws.Range("H22").Value = 1819.8889
$ws.Range("I22").Value = 1147.25
$ws.Range("K22").Value = 1147.25
$ws.Range("M22").Value = -852.25
# Row 27
$ws.Range("H27").Value = 1819.8889
$ws.Range("I27").Value = 1147.25
$ws.Range("K27").Value = 1147.25
$ws.Range("M27").Value = -1040.25
# Row 40
$ws.Range("H40").Value = 7318
$ws.Range("I40").Value = 7970
$ws.Range("K40").Value = 7970
$ws.Range("M40").Value = -7834
# Row 68
$ws.Range("H68").Value = 9950.1
$ws.Range("I68").Value = 12600.134
$ws.Range("K68").Value = 12600.134
$ws.Range("M68").Value = -11851.134
# Row 71
$ws.Range("H71").Value = 9950.1
$ws.Range("I71").Value = 12600.134
$ws.Range("K71").Value = 63000.67
$ws.Range("M71").Value = -59256.67
# Row 132
$ws.Range("H132").Value = 5008.6963
$ws.Range("I132").Value = 5067.3076
$ws.Range("J132").Value = 4874.2354
$ws.Range("K132").Value = 15201.9228
$ws.Range("L132").Value = 14622.7062
$ws.Range("M132").Value = -12671.9228
$ws.Range("N132").Value = -19682.7062

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
# Row 105
$ws.Range("H105").Value = 22352.4
$ws.Range("J105").Value = 22352.4
$ws.Range("L105").Value = 22352.4
$ws.Range("N105").Value = -29340.4

